$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a "history" data row (columns A-G) on a worksheet, matching
# the layout used throughout these sheets:
#   A = Run Date    (text that looks like a date, e.g. "2023-04-19")
#   B = Run Time    (numeric date/time serial, formatted as a date-time)
#   C = Sprint Name (text)
#   D = Total Cases (number)
#   E = Pass Cases  (number)
#   F = Fail Cases  (number)
#   G = Time Taken  (number)
#
# Parameters are positional (named parameter binding is unreliable here):
#   1 ws           worksheet COM object
#   2 row           destination row number
#   3 runDate       text value for column A
#   4 runTime       numeric value for column B
#   5 sprintName    text value for column C
#   6 totalCases    numeric value for column D
#   7 passCases     numeric value for column E
#   8 failCases     numeric value for column F
#   9 timeTaken     numeric value for column G
#  10 templateRow   (optional, 0 = none) row number to copy formatting from
# ---------------------------------------------------------------------------
function Write-HistoryRow {
    param($ws, $row, $runDate, $runTime, $sprintName, $totalCases, $passCases, $failCases, $timeTaken, $templateRow)

    # Column A - plain text that happens to look like a date. Force text
    # storage (NumberFormat "@") before assignment so the engine does not
    # auto-detect/convert it into a real date serial number, then drop the
    # temporary "@" number format back to a plain/general look.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $runDate
    $ws.Cells.Item($row, 1).Style = "Normal"

    # Column B - actual numeric date/time serial, shown with a date-time format.
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $runTime

    # Column C - plain text.
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $sprintName
    $ws.Cells.Item($row, 3).Style = "Normal"

    # Columns D-G - plain numbers.
    $ws.Cells.Item($row, 4).Value = $totalCases
    $ws.Cells.Item($row, 5).Value = $passCases
    $ws.Cells.Item($row, 6).Value = $failCases
    $ws.Cells.Item($row, 7).Value = $timeTaken

    # Pick up the same look/formatting as a neighboring existing data row
    # (general formatting for most columns, date-time formatting for B)
    # without disturbing the values we just wrote.
    if ($templateRow -gt 0) {
        $ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
        $ws.Range("A" + $row + ":G" + $row).PasteSpecial(-4122)
    }
}

# ---------------------------------------------------------------------------
# Sheet "AMSIN" (new tenant run-history rows for 2023-04-19 and 2023-04-20,
# plus a small correction to the existing 2023-04-18 row's run time / style)
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Row 67 already exists - only the run time (B67) actually changes value;
# re-apply the standard row formatting (matching row 66) while preserving
# all the existing cell values.
$wsAmsin.Cells.Item(67, 2).Value = 45034.63282362268
$wsAmsin.Range("A66:G66").Copy()
$wsAmsin.Range("A67:G67").PasteSpecial(-4122)

Write-HistoryRow $wsAmsin 68 "2023-04-19" 45035.69645784723 "176scndcyc" 155 154 1 4.1 66
Write-HistoryRow $wsAmsin 69 "2023-04-20" 45036.41314645833 "176fnlruntest" 155 154 1 3.49 68

# ---------------------------------------------------------------------------
# Sheet "BETA" (new tenant run-history row for 2023-04-20)
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Write-HistoryRow $wsBeta 35 "2023-04-20" 45036.5189090162 "176beta" 155 155 0 3.19 34

# ---------------------------------------------------------------------------
# Sheet "AMS" (new tenant run-history row for 2023-05-08)
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

Write-HistoryRow $wsAms 48 "2023-05-08" 45054.54037847627 "176htfxtrl" 155 155 0 3.02 0

"Done"
